# Update the cryptos list with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.310.16"
$ws.Range("E2").Value = "  +0.17%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.869.53"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "235.15"
$ws.Range("E5").Value = "  -0.97%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.05%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.33%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2870"
$ws.Range("E8").Value = "  +0.23%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06580"
$ws.Range("E9").Value = "  +0.55%  "

# Row 10 - Solana
$ws.Range("D10").Value = "21.64"
$ws.Range("E10").Value = "  -2.42%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.08020"
$ws.Range("E11").Value = "  +1.59%  "

# Row 12 - Litecoin
$ws.Range("D12").Value = "96.95"
$ws.Range("E12").Value = "  -0.86%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.871.31"
$ws.Range("E13").Value = "  +0.38%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.39%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.6853"
$ws.Range("E15").Value = "  +0.76%  "

# Row 16 - BitcoinCash
$ws.Range("D16").Value = "269.03"
$ws.Range("E16").Value = "  -3.26%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "30.319.79"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18 - Avalanche
$ws.Range("D18").Value = "14.05"
$ws.Range("E18").Value = "  +3.41%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.000007617"
$ws.Range("E19").Value = "  +3.76%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.10%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.116.32"
$ws.Range("E21").Value = "  +0.26%  "

# Row 22 - BinanceUSD
$ws.Range("E22").Value = "  +0.09%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.258"
$ws.Range("E23").Value = "  -2.31%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "6.217"
$ws.Range("E24").Value = "  +0.43%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "9.386"
$ws.Range("E25").Value = "  +1.17%  "

# Row 26 - Monero
$ws.Range("D26").Value = "167.66"
$ws.Range("E26").Value = "  -0.46%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -0.89%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "1.951"
$ws.Range("E28").Value = "  +0.44%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -1.21%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "0.09883"
$ws.Range("E30").Value = "  +0.75%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "4.362"
$ws.Range("E31").Value = "  -0.47%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.09%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.04712"
$ws.Range("E34").Value = "  -0.79%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  -0.48%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.7002"
$ws.Range("E36").Value = "  -0.75%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.32%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01884"
$ws.Range("E38").Value = "  +0.43%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "2.701"
$ws.Range("E39").Value = "  +3.04%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "6.266"
$ws.Range("E40").Value = "  -0.43%  "

# Row 41 - Aave
$ws.Range("D41").Value = "71.94"
$ws.Range("E41").Value = "  -5.85%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "1.956"
$ws.Range("E42").Value = "  +0.01%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "0.8415"
$ws.Range("E43").Value = "  -1.10%  "

# Rows 44/45 - PaxDollar and TheSandbox swap order
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4165"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46 - Quant
$ws.Range("D46").Value = "102.85"
$ws.Range("E46").Value = "  -0.23%  "

# Row 47 - Aptos
$ws.Range("D47").Value = "7.061"
$ws.Range("E47").Value = "  -2.18%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "9.089"
$ws.Range("E48").Value = "  -2.39%  "

# Row 49 - Maker
$ws.Range("D49").Value = "910.21"
$ws.Range("E49").Value = "  -4.09%  "

# Row 50 - Elrond
$ws.Range("E50").Value = "  +0.56%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.05704"
